$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 86/87, pushing the existing rows 86-117 down to 88-119.
$ws.Rows("86:87").Insert()

# Row 86: new weekly price observation.
$ws.Range("A86").Value = 9
$ws.Range("B86").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C86").Value = "Metropolitana"
$ws.Range("D86").Value = 44900
$ws.Range("E86").Value = 13
$ws.Range("F86").Value = 100114002
$ws.Range("G86").Value = "Camote"
$ws.Range("H86").Value = "Sin especificar"
$ws.Range("I86").Value = "Primera"
$ws.Range("J86").Value = 790
$ws.Range("K86").Value = 16000
$ws.Range("L86").Value = 18000
$ws.Range("M86").Value = 16987
$ws.Range("N86").Value = "$/caja 18 kilos"
$ws.Range("O86").Value = "Perú"
$ws.Range("P86").Value = 944
$ws.Range("Q86").Value = 18
$ws.Range("R86").Value = "Hortaliza"

# Row 87: new weekly price observation.
$ws.Range("A87").Value = 9
$ws.Range("B87").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C87").Value = "Metropolitana"
$ws.Range("D87").Value = 44900
$ws.Range("E87").Value = 13
$ws.Range("F87").Value = 100114002
$ws.Range("G87").Value = "Camote"
$ws.Range("H87").Value = "Sin especificar"
$ws.Range("I87").Value = "Primera"
$ws.Range("J87").Value = 520
$ws.Range("K87").Value = 14000
$ws.Range("L87").Value = 16000
$ws.Range("M87").Value = 15000
$ws.Range("N87").Value = "$/malla 18 kilos"
$ws.Range("O87").Value = "Perú"
$ws.Range("P87").Value = 833
$ws.Range("Q87").Value = 18
$ws.Range("R87").Value = "Hortaliza"
